$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits right
#    after the H1 title ("Play Cubes 2 Free and Enjoy a Unique Slot Game").
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
if ($metaPara.Range.Text -like "Meta description*") {
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# Helper XML fragments (OOXML "pkg:package" wrapper accepted by InsertXML).
# ---------------------------------------------------------------------------
$xmlNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function New-WordXmlPackage([string]$paragraphXml) {
    return "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" +
           "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" +
           "<pkg:xmlData><w:document $xmlNs><w:body>$paragraphXml</w:body></w:document></pkg:xmlData>" +
           "</pkg:part></pkg:package>"
}

# ---------------------------------------------------------------------------
# 2) Insert a new bold paragraph ("Play Cubes 2 Free and Enjoy a Unique Slot
#    Game") right before the closing "Create a cartoon-style..." paragraph.
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$lastPara.Range.InsertParagraphBefore()

$boldParaXml = "<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Cubes 2 Free and Enjoy a Unique Slot Game</w:t></w:r></w:p>"
$newBoldPara = $d.Paragraphs.Item($n)
$newBoldPara.Range.InsertXML((New-WordXmlPackage $boldParaXml))

# ---------------------------------------------------------------------------
# 3) Replace the text of the final paragraph (still italic) with the new
#    meta-description-style sentence. This paragraph is the very last one in
#    the document body, so its paragraph mark cannot simply be deleted; we
#    temporarily add a throw-away paragraph after it, rewrite this paragraph
#    in place, then drop the now-empty spare paragraph.
# ---------------------------------------------------------------------------
$n2 = $d.Paragraphs.Count
$finalPara = $d.Paragraphs.Item($n2)
$finalPara.Range.InsertParagraphAfter()

$italicParaXml = "<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Experience the excitement of Cubes 2, a high-volatility slot game based on Rubik&#39;s cube. Play for free now!</w:t></w:r></w:p>"
$targetPara = $d.Paragraphs.Item($n2)
$targetPara.Range.InsertXML((New-WordXmlPackage $italicParaXml))

$spareIndex = $d.Paragraphs.Count
$d.Paragraphs.Item($spareIndex).Range.Delete()

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
